$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a trailing period to the three assumption cells that were missing one.
$ws.Range("E5").Value = "Chose 2019 GFDx number (33 g/gday) over 2000-2009 estimate of wheat (not specifically flour) (35 g/day) and over 1997-2000 survey estimate (17 g/day)."
$ws.Range("E6").Value = "Chose 2017 GFDx estimate (126 g/day) over number from 2007 paper (67 g/day) that estimated from FAO balance sheets, and a 1997-2000 report that estimated from FAO balance sheets."
$ws.Range("E8").Value = "Chose 2017 GFDx estimte (283 g/day) over a 2007 paper that estimated (288 g/day) from FAO balance sheets and a 1997-2000 paper that estimated (204 g/day) from FAO balance sheets."

# Adjust row heights that changed.
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 60
$ws.Rows.Item(5).RowHeight = 30

# Update the view's active selection (matches the final cursor position).
$ws.Range("E9").Select()
